$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 478 (2021-06-17): the new-cases count for that day is corrected from 0 to 7.
$ws.Range("C478").Value = 7

# L478/M478 were holding the text placeholder "0" (the column is formatted as
# Text, numFmtId 49 / "@"). Re-enter them as real numbers while keeping the
# column's Text display format: drop to General long enough to write a true
# numeric value, then restore "@" so the style id collapses back to the
# original text-format style.
$ws.Range("L478").ClearContents()
$ws.Range("L478").NumberFormat = "General"
$ws.Range("L478").Value = 0
$ws.Range("L478").NumberFormat = "@"

$ws.Range("M478").ClearContents()
$ws.Range("M478").NumberFormat = "General"
$ws.Range("M478").Value = 0
$ws.Range("M478").NumberFormat = "@"

# Row 479 (2021-06-18): data for that day is now filled in - 2 new cases.
$ws.Range("C479").Value = 2

$ws.Range("L479").ClearContents()
$ws.Range("L479").NumberFormat = "General"
$ws.Range("L479").Value = 0
$ws.Range("L479").NumberFormat = "@"

$ws.Range("M479").ClearContents()
$ws.Range("M479").NumberFormat = "General"
$ws.Range("M479").Value = 0
$ws.Range("M479").NumberFormat = "@"

# Move the frozen-pane viewport down a couple of rows and update the
# remembered selection to where the editor was working.
$win = $excel.ActiveWindow
$win.ScrollRow = 468
$win.ScrollColumn = 2
$ws.Range("Q482").Select()
